$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated price / volume(1h) values from the latest crypto data pull.
# Force text storage (NumberFormat "@") then ClearFormats so the cell keeps
# its original (default) style -- only the text content changes, matching
# the source data which stores these as plain strings.

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '27.512.39'
$c.ClearFormats()

$c = $ws.Range('E2')
$c.NumberFormat = "@"
$c.Value = '  +0.44%  '
$c.ClearFormats()

$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '1.742.24'
$c.ClearFormats()

$c = $ws.Range('E3')
$c.NumberFormat = "@"
$c.Value = '  -0.31%  '
$c.ClearFormats()

$c = $ws.Range('E4')
$c.NumberFormat = "@"
$c.Value = '  -0.07%  '
$c.ClearFormats()

$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '322.97'
$c.ClearFormats()

$c = $ws.Range('E5')
$c.NumberFormat = "@"
$c.Value = '  +0.21%  '
$c.ClearFormats()

$c = $ws.Range('E6')
$c.NumberFormat = "@"
$c.Value = '  +0.01%  '
$c.ClearFormats()

$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '0.4480'
$c.ClearFormats()

$c = $ws.Range('E7')
$c.NumberFormat = "@"
$c.Value = '  +5.69%  '
$c.ClearFormats()

$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '0.3524'
$c.ClearFormats()

$c = $ws.Range('E8')
$c.NumberFormat = "@"
$c.Value = '  -2.05%  '
$c.ClearFormats()

$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.07374'
$c.ClearFormats()

$c = $ws.Range('E9')
$c.NumberFormat = "@"
$c.Value = '  -1.52%  '
$c.ClearFormats()

$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '41.34'
$c.ClearFormats()

$c = $ws.Range('E10')
$c.NumberFormat = "@"
$c.Value = '  -1.73%  '
$c.ClearFormats()

$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '1.076'
$c.ClearFormats()

$c = $ws.Range('E11')
$c.NumberFormat = "@"
$c.Value = '  -1.95%  '
$c.ClearFormats()

$c = $ws.Range('E12')
$c.NumberFormat = "@"
$c.Value = '  -0.13%  '
$c.ClearFormats()

$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '20.43'
$c.ClearFormats()

$c = $ws.Range('E13')
$c.NumberFormat = "@"
$c.Value = '  -1.03%  '
$c.ClearFormats()

$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '5.899'
$c.ClearFormats()

$c = $ws.Range('E14')
$c.NumberFormat = "@"
$c.Value = '  -2.05%  '
$c.ClearFormats()

$c = $ws.Range('E15')
$c.NumberFormat = "@"
$c.Value = '  -2.09%  '
$c.ClearFormats()

$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '1.742.10'
$c.ClearFormats()

$c = $ws.Range('E16')
$c.NumberFormat = "@"
$c.Value = '  -0.37%  '
$c.ClearFormats()

$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '91.50'
$c.ClearFormats()

$c = $ws.Range('E17')
$c.NumberFormat = "@"
$c.Value = '  -1.27%  '
$c.ClearFormats()

$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '0.00001052'
$c.ClearFormats()

$c = $ws.Range('E18')
$c.NumberFormat = "@"
$c.Value = '  -1.51%  '
$c.ClearFormats()

$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '0.06363'
$c.ClearFormats()

$c = $ws.Range('E19')
$c.NumberFormat = "@"
$c.Value = '  -0.13%  '
$c.ClearFormats()

$c = $ws.Range('E20')
$c.NumberFormat = "@"
$c.Value = '  +0.04%  '
$c.ClearFormats()

$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '16.80'
$c.ClearFormats()

$c = $ws.Range('E21')
$c.NumberFormat = "@"
$c.Value = '  -1.40%  '
$c.ClearFormats()

$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '5.721'
$c.ClearFormats()

$c = $ws.Range('E22')
$c.NumberFormat = "@"
$c.Value = '  -2.72%  '
$c.ClearFormats()

$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '27.548.71'
$c.ClearFormats()

$c = $ws.Range('E23')
$c.NumberFormat = "@"
$c.Value = '  +0.32%  '
$c.ClearFormats()

$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '11.11'
$c.ClearFormats()

$c = $ws.Range('E24')
$c.NumberFormat = "@"
$c.Value = '  -0.74%  '
$c.ClearFormats()

$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '2.098'
$c.ClearFormats()

$c = $ws.Range('E25')
$c.NumberFormat = "@"
$c.Value = '  +0.35%  '
$c.ClearFormats()

$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '161.81'
$c.ClearFormats()

$c = $ws.Range('E26')
$c.NumberFormat = "@"
$c.Value = '  +0.03%  '
$c.ClearFormats()

$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '20.03'
$c.ClearFormats()

$c = $ws.Range('E27')
$c.NumberFormat = "@"
$c.Value = '  -1.15%  '
$c.ClearFormats()

$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '1.944.02'
$c.ClearFormats()

$c = $ws.Range('E28')
$c.NumberFormat = "@"
$c.Value = '  -0.19%  '
$c.ClearFormats()

$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '124.67'
$c.ClearFormats()

$c = $ws.Range('E29')
$c.NumberFormat = "@"
$c.Value = '  +0.70%  '
$c.ClearFormats()

$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '2.025'
$c.ClearFormats()

$c = $ws.Range('E31')
$c.NumberFormat = "@"
$c.Value = '  -5.17%  '
$c.ClearFormats()

$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '0.09034'
$c.ClearFormats()

$c = $ws.Range('E32')
$c.NumberFormat = "@"
$c.Value = '  +1.86%  '
$c.ClearFormats()

$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '3.652'
$c.ClearFormats()

$c = $ws.Range('E33')
$c.NumberFormat = "@"
$c.Value = '  +0.14%  '
$c.ClearFormats()

$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '5.364'
$c.ClearFormats()

$c = $ws.Range('E34')
$c.NumberFormat = "@"
$c.Value = '  -2.93%  '
$c.ClearFormats()

$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '0.02268'
$c.ClearFormats()

$c = $ws.Range('E35')
$c.NumberFormat = "@"
$c.Value = '  -0.57%  '
$c.ClearFormats()

$c = $ws.Range('E36')
$c.NumberFormat = "@"
$c.Value = '  -4.74%  '
$c.ClearFormats()

$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '0.05987'
$c.ClearFormats()

$c = $ws.Range('E37')
$c.NumberFormat = "@"
$c.Value = '  -0.17%  '
$c.ClearFormats()

$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '0.2058'
$c.ClearFormats()

$c = $ws.Range('E38')
$c.NumberFormat = "@"
$c.Value = '  -1.86%  '
$c.ClearFormats()

$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '0.6239'
$c.ClearFormats()

$c = $ws.Range('E39')
$c.NumberFormat = "@"
$c.Value = '  -1.42%  '
$c.ClearFormats()

$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '4.873'
$c.ClearFormats()

$c = $ws.Range('E40')
$c.NumberFormat = "@"
$c.Value = '  -1.34%  '
$c.ClearFormats()

$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '1.184'
$c.ClearFormats()

$c = $ws.Range('E41')
$c.NumberFormat = "@"
$c.Value = '  +0.11%  '
$c.ClearFormats()

$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '1.377'
$c.ClearFormats()

$c = $ws.Range('E42')
$c.NumberFormat = "@"
$c.Value = '  -0.63%  '
$c.ClearFormats()

$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '7.702'
$c.ClearFormats()

$c = $ws.Range('E43')
$c.NumberFormat = "@"
$c.Value = '  -2.08%  '
$c.ClearFormats()

$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '13.09'
$c.ClearFormats()

$c = $ws.Range('E44')
$c.NumberFormat = "@"
$c.Value = '  -2.08%  '
$c.ClearFormats()

$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '3.702'
$c.ClearFormats()

$c = $ws.Range('E45')
$c.NumberFormat = "@"
$c.Value = '  +0.30%  '
$c.ClearFormats()

$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '0.5789'
$c.ClearFormats()

$c = $ws.Range('E46')
$c.NumberFormat = "@"
$c.Value = '  -1.38%  '
$c.ClearFormats()

$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '121.90'
$c.ClearFormats()

$c = $ws.Range('E47')
$c.NumberFormat = "@"
$c.Value = '  -0.33%  '
$c.ClearFormats()

$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '1.920'
$c.ClearFormats()

$c = $ws.Range('E48')
$c.NumberFormat = "@"
$c.Value = '  -2.44%  '
$c.ClearFormats()

$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '0.06836'
$c.ClearFormats()

$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '1.111'
$c.ClearFormats()

$c = $ws.Range('E50')
$c.NumberFormat = "@"
$c.Value = '  -4.77%  '
$c.ClearFormats()

$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '71.15'
$c.ClearFormats()

$c = $ws.Range('E51')
$c.NumberFormat = "@"
$c.Value = '  -2.61%  '
$c.ClearFormats()

Write-Output "Applied 91 cell updates (Price/Volume(1h) columns)"
